# Scheduled-runner refresh of market-price-derived columns (H:N) across
# several leve-profit tables (one per crafting job sheet). Only numeric
# input cells change; no formulas, styles, or structure are touched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 669.7
$ws.Range("I19").Value = 459.8
$ws.Range("J19").Value = 879.6
$ws.Range("K19").Value = 459.8
$ws.Range("L19").Value = 879.6
$ws.Range("M19").Value = -284.8
$ws.Range("N19").Value = -1229.6

$ws.Range("H98").Value = 1874.625
$ws.Range("I98").Value = 1828.7354
$ws.Range("J98").Value = 2134.6667
$ws.Range("K98").Value = 1828.7354
$ws.Range("L98").Value = 2134.6667
$ws.Range("M98").Value = -330.7354
$ws.Range("N98").Value = -5130.6667

$ws.Range("H122").Value = 1874.625
$ws.Range("I122").Value = 1828.7354
$ws.Range("J122").Value = 2134.6667
$ws.Range("K122").Value = 5486.206200000001
$ws.Range("L122").Value = 6404.000100000001
$ws.Range("M122").Value = -3036.206200000001
$ws.Range("N122").Value = -11304.0001

$ws.Range("H132").Value = 25006742
$ws.Range("I132").Value = 37042932
$ws.Range("J132").Value = 8499.615
$ws.Range("K132").Value = 111128796
$ws.Range("L132").Value = 25498.845
$ws.Range("M132").Value = -111126266
$ws.Range("N132").Value = -30558.845

$ws.Range("H137").Value = 2967.2693
$ws.Range("I137").Value = 2979.7
$ws.Range("J137").Value = 2959.5
$ws.Range("K137").Value = 8939.099999999999
$ws.Range("L137").Value = 8878.5
$ws.Range("M137").Value = -6389.099999999999
$ws.Range("N137").Value = -13978.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2451.9011
$ws.Range("I32").Value = 2417.147
$ws.Range("J32").Value = 2633.6924
$ws.Range("K32").Value = 2417.147
$ws.Range("L32").Value = 2633.6924
$ws.Range("M32").Value = -2130.147
$ws.Range("N32").Value = -3207.6924

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 806.5161000000001
$ws.Range("I20").Value = 800.0526
$ws.Range("J20").Value = 816.75
$ws.Range("K20").Value = 800.0526
$ws.Range("L20").Value = 816.75
$ws.Range("M20").Value = -553.0526
$ws.Range("N20").Value = -1310.75

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H64").Value = 1703
$ws.Range("J64").Value = 1955
$ws.Range("L64").Value = 1955
$ws.Range("N64").Value = -2405

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H67").Value = 1703
$ws.Range("J67").Value = 1955
$ws.Range("L67").Value = 1955
$ws.Range("N67").Value = -3515

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2456.3794
$ws.Range("I31").Value = 2414.6155
$ws.Range("J31").Value = 2468.4443
$ws.Range("K31").Value = 2414.6155
$ws.Range("L31").Value = 2468.4443
$ws.Range("M31").Value = -2119.6155
$ws.Range("N31").Value = -3058.4443

$ws.Range("H34").Value = 2456.3794
$ws.Range("I34").Value = 2414.6155
$ws.Range("J34").Value = 2468.4443
$ws.Range("K34").Value = 2414.6155
$ws.Range("L34").Value = 2468.4443
$ws.Range("M34").Value = -2212.6155
$ws.Range("N34").Value = -2872.4443

$ws.Range("H109").Value = 49450
$ws.Range("J109").Value = 49450
$ws.Range("L109").Value = 49450
$ws.Range("N109").Value = -51530

$ws.Range("H122").Value = 3942.5789
$ws.Range("I122").Value = 3181.9375
$ws.Range("J122").Value = 7999.3335
$ws.Range("K122").Value = 9545.8125
$ws.Range("L122").Value = 23998.0005
$ws.Range("M122").Value = -7095.8125
$ws.Range("N122").Value = -28898.0005

$ws.Range("H134").Value = 1437.4688
$ws.Range("I134").Value = 1443.1111
$ws.Range("J134").Value = 1430.2142
$ws.Range("K134").Value = 4329.3333
$ws.Range("L134").Value = 4290.642599999999
$ws.Range("M134").Value = -1794.3333
$ws.Range("N134").Value = -9360.642599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6459.4736
$ws.Range("I5").Value = 1079.4
$ws.Range("J5").Value = 8380.929
$ws.Range("K5").Value = 3238.2
$ws.Range("L5").Value = 25142.787
$ws.Range("M5").Value = -3126.2
$ws.Range("N5").Value = -25366.787

$ws.Range("H68").Value = 6655.143
$ws.Range("I68").Value = 2416.8333
$ws.Range("J68").Value = 7811.0454
$ws.Range("K68").Value = 7250.499899999999
$ws.Range("L68").Value = 23433.1362
$ws.Range("M68").Value = -6439.499899999999
$ws.Range("N68").Value = -25055.1362

$ws.Range("H71").Value = 6655.143
$ws.Range("I71").Value = 2416.8333
$ws.Range("J71").Value = 7811.0454
$ws.Range("K71").Value = 21751.4997
$ws.Range("L71").Value = 70299.4086
$ws.Range("M71").Value = -17695.4997
$ws.Range("N71").Value = -78411.4086

$ws.Range("H132").Value = 2118.5386
$ws.Range("I132").Value = 1514
$ws.Range("J132").Value = 2299.9
$ws.Range("K132").Value = 13626
$ws.Range("L132").Value = 20699.1
$ws.Range("M132").Value = -11096
$ws.Range("N132").Value = -25759.1

$ws.Range("H135").Value = 6459.4736
$ws.Range("I135").Value = 1079.4
$ws.Range("J135").Value = 8380.929
$ws.Range("K135").Value = 9714.6
$ws.Range("L135").Value = 75428.361
$ws.Range("M135").Value = -7179.6
$ws.Range("N135").Value = -80498.361

$ws.Range("H137").Value = 24310198
$ws.Range("I137").Value = 20835904
$ws.Range("J137").Value = 27784490
$ws.Range("K137").Value = 62507712
$ws.Range("L137").Value = 83353470
$ws.Range("M137").Value = -62502612
$ws.Range("N137").Value = -83363670

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 25671008
$ws.Range("I122").Value = 37078430
$ws.Range("J122").Value = 4303.75
$ws.Range("K122").Value = 111235290
$ws.Range("L122").Value = 12911.25
$ws.Range("M122").Value = -111232840
$ws.Range("N122").Value = -17811.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1028.3462
$ws.Range("I55").Value = 100.61539
$ws.Range("J55").Value = 1956.0769
$ws.Range("K55").Value = 100.61539
$ws.Range("L55").Value = 1956.0769
$ws.Range("M55").Value = 72.38461
$ws.Range("N55").Value = -2302.0769

$ws.Range("H122").Value = 3708.2163
$ws.Range("I122").Value = 3694.25
$ws.Range("J122").Value = 3797.6
$ws.Range("K122").Value = 11082.75
$ws.Range("L122").Value = 11392.8
$ws.Range("M122").Value = -8632.75
$ws.Range("N122").Value = -16292.8
